$d = $word.ActiveDocument

# Fix hyphenation artifact: "первом лич- ном году" -> "первом личном году"
$d.Content.Find.Execute("первом лич- ном году", $true, $false, $false, $false, $false,
                         $true, 1, $false, "первом личном году", 2)
